# Append: 2025-12-19 01:56 JST
# Update the "取得日時" (acquisition timestamp) column (A) for all data rows
# on the active sheet (ランサーズ) from the previous run timestamp to the
# new one, reflecting a fresh scrape pass at 2025-12-19 01:56:10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-12-19 01:56:10"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
